$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values regenerated to filter save games (columns B-E, with G = sum of B:E)
$data = @{
    2 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 6.15379541431027 }
    3 = @{ B = 0.2881169905109251; C = 9.983522426115931;  D = 3.223369029078222;  E = 13.86384647080068;   G = 27.35885491650576 }
    4 = @{ B = 0.2881169905109251; C = 1.626987699542094;  D = 0.1496068669990043; E = 0.5333859586016987;  G = 2.598097515653722 }
    5 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 0.1496068669990043; E = 13.86384647080068;   G = 18.91276827552123 }
    6 = @{ B = 0.6545652718822623; C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 3.536033448013082 }
    7 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 4.327115817150455 }
    8 = @{ B = 1.445647641019636;  C = 1.626987699542094;  D = 3.223369029078222;  E = 0.5333859586016987;  G = 6.82939032824165 }
    9 = @{ B = 3.272327238179451;  C = 1.626987699542094;  D = 3.223369029078222;  E = 0.5333859586016987;  G = 8.656069925401464 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
